$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column B
$ws.Range("B2").Value = 124
$ws.Range("B3").Value = 99

# Add new rows 4 and 5
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 42
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 33

# Copy the formatting (border/bold/centered style) from A3 to the new A-column cells
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0
